## "Prouct excel sheet added"
## Adds a new "product" worksheet (after the existing "Login" sheet) containing
## four product blocks (header row + one data row each), highlights every
## header row with a yellow fill, sizes a few columns, and leaves the new
## sheet as the active tab/selection.

$wb  = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item(1)

# Insert the new sheet right after "Login" -- this also makes it the active
# sheet (tabSelected on the sheet / activeTab on the workbook).
$ws = $wb.Worksheets.Add($null, $login)
$ws.Name = "product"

# ---- Header rows (productname / price / commissionrate / unitqty / stockquantity) ----
# Style first, then values, so the header cells keep t="s" + s="1" together.
$ws.Range("A1:E1").Interior.Color = 65535
$ws.Range("A4:E4").Interior.Color = 65535
$ws.Range("A7:E7").Interior.Color = 65535
$ws.Range("A10:E10").Interior.Color = 65535

$headers = "productname","price","commissionrate","unitqty","stockquantity"
$headerRows = 1,4,7,10
foreach ($r in $headerRows) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $headers[$c - 1]
    }
}

# ---- Product data rows ----
# Written in this order so the shared-string table fills up the same way it
# did in the source workbook (Ikegai, Mindset, Rich Dad poor Dad, Atomic
# Habbits) even though the rows don't land in that order on the sheet.
$ws.Cells.Item(5, 1).Value = "Ikegai"
$ws.Cells.Item(8, 1).Value = "Mindset"
$ws.Cells.Item(11, 1).Value = "Rich Dad poor Dad"
$ws.Cells.Item(2, 1).Value = "Atomic Habbits"

# price, commissionrate, unitqty, stockquantity for each product row
$ws.Cells.Item(2, 2).Value = 150
$ws.Cells.Item(2, 3).Value = 15
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 25

$ws.Cells.Item(5, 2).Value = 120
$ws.Cells.Item(5, 3).Value = 20
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 30

$ws.Cells.Item(8, 2).Value = 100
$ws.Cells.Item(8, 3).Value = 15
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 15

$ws.Cells.Item(11, 2).Value = 160
$ws.Cells.Item(11, 3).Value = 10
$ws.Cells.Item(11, 4).Value = 3
$ws.Cells.Item(11, 5).Value = 40

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 17.43
$ws.Columns.Item(2).ColumnWidth = 14.43
$ws.Columns.Item(5).ColumnWidth = 13.43

# ---- Selection on the new sheet ----
$ws.Range("G8").Select() | Out-Null
